$d = $word.ActiveDocument
$w_ns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# Paragraph/run formatting shared by every line in this scene (sz/szCs 24, rtl 0,
# 1.15-line spacing) - matches the canonical OOXML shown in the diff.
$pPr = "<w:pPr><w:spacing w:line='276' w:lineRule='auto'/><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr></w:pPr>"
function RunXml($text) {
    return "<w:r><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/><w:rtl w:val='0'/></w:rPr><w:t xml:space='preserve'>$text</w:t></w:r>"
}
function ParaXml($text) {
    return "<w:p $w_ns>$pPr$(RunXml($text))</w:p>"
}

# Find the (single) paragraph whose text contains $needle and replace its whole
# contents with one run holding $newText, serialized the canonical way (single
# run, xml:space="preserve" present). Using InsertXML (rather than Find/Replace or
# Range.Text=) both lets us merge multiple runs into one AND forces the
# xml:space="preserve" flag that a plain text replace in this interop drops.
function ReplaceParaContaining($needle, $newText) {
    $paras = $d.Paragraphs
    $n = $paras.Count
    $i = 1
    $hit = 0
    while ($i -le $n) {
        $cand = $paras.Item($i)
        if ($cand.Range.Text.Contains($needle)) {
            $null = $cand.Range.InsertXML((ParaXml($newText)))
            $hit = 1
            $i = $n + 1
        } else {
            $i = $i + 1
        }
    }
    return $hit
}

$ellipsis = [char]8230
$rsquo = [char]8217

# 1) "Petra (neutral embarrassed): It was embarrassing..." - text itself is
#    unchanged; only the XML serialization noise (missing xml:space="preserve")
#    needs to be normalized away.
$text1 = "Petra (neutral embarrassed): It was embarrassing$ellipsis"
if (-not (ReplaceParaContaining "neutral embarrassed): It was embarrassing" $text1)) {
    throw "Could not find the 'Petra (neutral embarrassed)' paragraph"
}

# 2) Merge the 3 runs "Teacher (" / "neutral disappointed" / "): You know..." into
#    a single run.
$text2 = "Teacher (neutral disappointed): You know$ellipsis"
if (-not (ReplaceParaContaining "neutral disappointed" $text2)) {
    throw "Could not find the 'Teacher (neutral disappointed)' paragraph"
}

# 3) Merge the 5 runs making up the "expressionless" line into a single run.
$text3 = "Teacher (neutral expressionless): But I get wanting to sleep in. If I could sleep in every day, I would$ellipsis"
if (-not (ReplaceParaContaining "expressionless): But I" $text3)) {
    throw "Could not find the 'Teacher (neutral expressionless)' paragraph"
}

# 4) "Teacher (arms_crossed smug): Oh, I heard..." -> "...curious): Oh, I heard..."
$text4 = "Teacher (arms_crossed curious): Oh, I heard from Asher that he" + $rsquo + "s dragging you along to a mixer or something?"
if (-not (ReplaceParaContaining "arms_crossed smug): Oh, I heard" $text4)) {
    throw "Could not find the 'Teacher (arms_crossed smug): Oh, I heard' paragraph"
}

# 5) "Teacher (neutral grinning): That's punishment..." -> "...arms_crossed smug): That's punishment..."
$text5 = "Teacher (arms_crossed smug): That" + $rsquo + "s punishment enough for you. Try not to make a fool of yourself."
if (-not (ReplaceParaContaining "neutral grinning): That" $text5)) {
    throw "Could not find the 'Teacher (neutral grinning)' paragraph"
}

# 6) "Teacher (neutral laughing): No problem." -> "...arms_crossed smiling_eyes_closed): No problem."
$text6 = "Teacher (arms_crossed smiling_eyes_closed): No problem."
if (-not (ReplaceParaContaining "neutral laughing): No problem" $text6)) {
    throw "Could not find the 'Teacher (neutral laughing)' paragraph"
}

# 7) "Teacher (neutral smiling): Now get back to class..." -> "...neutral neutral): Now get back to class..."
$text7 = "Teacher (neutral neutral): Now get back to class, I wanna eat."
if (-not (ReplaceParaContaining "neutral smiling): Now get back to class" $text7)) {
    throw "Could not find the 'Teacher (neutral smiling)' paragraph"
}

Write-Output "done"
